$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 3225.762
$ws_ALC.Range("I28").Value = 1144.6923
$ws_ALC.Range("J28").Value = 6607.5
$ws_ALC.Range("K28").Value = 1144.6923
$ws_ALC.Range("L28").Value = 6607.5
$ws_ALC.Range("M28").Value = -659.6922999999999
$ws_ALC.Range("N28").Value = -7577.5
$ws_ALC.Range("H88").Value = 3000
$ws_ALC.Range("J88").Value = 0
$ws_ALC.Range("L88").Value = 0
$ws_ALC.Range("N88").ClearContents()
$ws_ALC.Range("H91").Value = 3000
$ws_ALC.Range("J91").Value = 0
$ws_ALC.Range("L91").Value = 0
$ws_ALC.Range("N91").ClearContents()
$ws_ALC.Range("H96").Value = 176.09091
$ws_ALC.Range("J96").Value = 87.5
$ws_ALC.Range("L96").Value = 262.5
$ws_ALC.Range("N96").Value = -3008.5
$ws_ALC.Range("H100").Value = 1380.375
$ws_ALC.Range("I100").Value = 1529.4
$ws_ALC.Range("K100").Value = 1529.4
$ws_ALC.Range("M100").Value = -988.4000000000001
$ws_ALC.Range("H111").Value = 1101.2142
$ws_ALC.Range("I111").Value = 369.75
$ws_ALC.Range("K111").Value = 1109.25
$ws_ALC.Range("M111").Value = 1957.75
$ws_ALC.Range("H113").Value = 4004.8
$ws_ALC.Range("I113").Value = 4004.8
$ws_ALC.Range("K113").Value = 4004.8
$ws_ALC.Range("M113").Value = -750.8000000000002
$ws_ALC.Range("H137").Value = 4498.222
$ws_ALC.Range("I137").Value = 934.9167
$ws_ALC.Range("K137").Value = 2804.7501
$ws_ALC.Range("M137").Value = -254.7501000000002
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 1969.1428
$ws_ARM.Range("I61").Value = 932
$ws_ARM.Range("K61").Value = 932
$ws_ARM.Range("M61").Value = -720
$ws_ARM.Range("H63").Value = 5456.4287
$ws_ARM.Range("I63").Value = 2847.5
$ws_ARM.Range("K63").Value = 2847.5
$ws_ARM.Range("M63").Value = -2161.5
$ws_ARM.Range("H66").Value = 5456.4287
$ws_ARM.Range("I66").Value = 2847.5
$ws_ARM.Range("K66").Value = 14237.5
$ws_ARM.Range("M66").Value = -10805.5
$ws_ARM.Range("H110").Value = 1984.3334
$ws_ARM.Range("I110").Value = 1763.7142
$ws_ARM.Range("J110").Value = 2756.5
$ws_ARM.Range("K110").Value = 1763.7142
$ws_ARM.Range("L110").Value = 2756.5
$ws_ARM.Range("M110").Value = 281.2858000000001
$ws_ARM.Range("N110").Value = -6846.5
$ws_ARM.Range("H136").Value = 1969.1428
$ws_ARM.Range("I136").Value = 932
$ws_ARM.Range("K136").Value = 2796
$ws_ARM.Range("M136").Value = -246
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 823.1875
$ws_BSM.Range("I94").Value = 612.2857
$ws_BSM.Range("K94").Value = 612.2857
$ws_BSM.Range("M94").Value = -161.2857
$ws_BSM.Range("H134").Value = 1480.5454
$ws_BSM.Range("I134").Value = 1032.6
$ws_BSM.Range("J134").Value = 5960
$ws_BSM.Range("K134").Value = 3097.8
$ws_BSM.Range("L134").Value = 17880
$ws_BSM.Range("M134").Value = -562.7999999999997
$ws_BSM.Range("N134").Value = -22950
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H4").Value = 3127023.8
$ws_CRP.Range("J4").Value = 12502125
$ws_CRP.Range("L4").Value = 12502125
$ws_CRP.Range("N4").Value = -12502349
$ws_CRP.Range("H20").Value = 39375
$ws_CRP.Range("J20").Value = 39375
$ws_CRP.Range("L20").Value = 39375
$ws_CRP.Range("N20").Value = -39847
$ws_CRP.Range("H30").Value = 39375
$ws_CRP.Range("J30").Value = 39375
$ws_CRP.Range("L30").Value = 39375
$ws_CRP.Range("N30").Value = -39557
$ws_CRP.Range("H62").Value = 2899.75
$ws_CRP.Range("I62").Value = 2899.75
$ws_CRP.Range("J62").Value = 0
$ws_CRP.Range("K62").Value = 2899.75
$ws_CRP.Range("L62").Value = 0
$ws_CRP.Range("M62").ClearContents()
$ws_CRP.Range("N62").Value = -2275.75
$ws_CRP.Range("H65").Value = 2899.75
$ws_CRP.Range("I65").Value = 2899.75
$ws_CRP.Range("J65").Value = 0
$ws_CRP.Range("K65").Value = 14498.75
$ws_CRP.Range("L65").Value = 0
$ws_CRP.Range("M65").ClearContents()
$ws_CRP.Range("N65").Value = -11378.75
$ws_CRP.Range("H107").Value = 872.4211
$ws_CRP.Range("I107").Value = 350.6154
$ws_CRP.Range("K107").Value = 350.6154
$ws_CRP.Range("M107").Value = 1569.3846
$ws_CRP.Range("H122").Value = 985.1539
$ws_CRP.Range("I122").Value = 985.1539
$ws_CRP.Range("K122").Value = 2955.4617
$ws_CRP.Range("M122").Value = -505.4616999999998
$ws_CRP.Range("H128").Value = 39375
$ws_CRP.Range("J128").Value = 39375
$ws_CRP.Range("L128").Value = 39375
$ws_CRP.Range("N128").Value = -49335
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H14").Value = 641.2
$ws_CUL.Range("I14").Value = 641.2
$ws_CUL.Range("K14").Value = 1923.6
$ws_CUL.Range("M14").Value = -1750.6
$ws_CUL.Range("H68").Value = 0
$ws_CUL.Range("I68").Value = 0
$ws_CUL.Range("J68").Value = 0
$ws_CUL.Range("K68").Value = 0
$ws_CUL.Range("L68").ClearContents()
$ws_CUL.Range("M68").ClearContents()
$ws_CUL.Range("N68").Value = 0
$ws_CUL.Range("H71").Value = 0
$ws_CUL.Range("I71").Value = 0
$ws_CUL.Range("J71").Value = 0
$ws_CUL.Range("K71").Value = 0
$ws_CUL.Range("L71").ClearContents()
$ws_CUL.Range("M71").ClearContents()
$ws_CUL.Range("N71").Value = 0
$ws_CUL.Range("H92").Value = 765.1
$ws_CUL.Range("I92").Value = 643.875
$ws_CUL.Range("K92").Value = 1931.625
$ws_CUL.Range("M92").Value = -683.625
$ws_CUL.Range("H123").Value = 0
$ws_CUL.Range("I123").Value = 0
$ws_CUL.Range("K123").Value = 0
$ws_CUL.Range("M123").ClearContents()
$ws_CUL.Range("H131").Value = 1712.2858
$ws_CUL.Range("I131").Value = 1199.4
$ws_CUL.Range("J131").Value = 2994.5
$ws_CUL.Range("K131").Value = 3598.2
$ws_CUL.Range("L131").Value = 8983.5
$ws_CUL.Range("M131").Value = 1441.8
$ws_CUL.Range("N131").Value = -19063.5
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H116").Value = 0
$ws_GSM.Range("J116").Value = 0
$ws_GSM.Range("L116").ClearContents()
$ws_GSM.Range("N116").Value = 0
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H16").Value = 0
$ws_WVR.Range("I16").Value = 0
$ws_WVR.Range("K16").Value = 0
$ws_WVR.Range("M16").ClearContents()
$ws_WVR.Range("H107").Value = 763.7692
$ws_WVR.Range("I107").Value = 769.44446
$ws_WVR.Range("K107").Value = 2308.33338
$ws_WVR.Range("M107").Value = -388.33338
$ws_WVR.Range("H122").Value = 1434
$ws_WVR.Range("I122").Value = 1249.75
$ws_WVR.Range("K122").Value = 3749.25
$ws_WVR.Range("M122").Value = -1299.25
$ws_WVR.Range("H136").Value = 3764.7407
$ws_WVR.Range("I136").Value = 2985.5293
$ws_WVR.Range("K136").Value = 8956.5879
$ws_WVR.Range("M136").Value = -6406.5879
